# Add DT line 2 MRT stations, and fix the corrupted Woodlands Chinese name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix corrupted Chinese station name for Woodlands (NS9) in column C, row 40.
$ws.Cells.Item(40, 3).Value = "兀兰站"

# Append the Downtown Line (2nd stage) stations: DT1 .. DT13
# Columns: A=Number, B=Name(en), C=Name(zh), D=Latitude, E=Longitude, F=ZIP, G=Street

$ws.Cells.Item(144, 1).Value = "DT1"
$ws.Cells.Item(144, 2).Value = "Bukit Panjang"
$ws.Cells.Item(144, 3).Value = "武吉班让"
$ws.Cells.Item(144, 4).Value = 1.378197
$ws.Cells.Item(144, 5).Value = 103.763545
$ws.Cells.Item(144, 6).Value = 678270
$ws.Cells.Item(144, 7).Value = "15 Petir Rd"

$ws.Cells.Item(145, 1).Value = "DT2"
$ws.Cells.Item(145, 2).Value = "Cashew"
$ws.Cells.Item(145, 4).Value = 1.368975
$ws.Cells.Item(145, 5).Value = 103.764803

$ws.Cells.Item(146, 1).Value = "DT3"
$ws.Cells.Item(146, 2).Value = "Hillview"
$ws.Cells.Item(146, 4).Value = 1.362472
$ws.Cells.Item(146, 5).Value = 103.767389

$ws.Cells.Item(147, 1).Value = "DT5"
$ws.Cells.Item(147, 2).Value = "Beauty World"
$ws.Cells.Item(147, 4).Value = 1.341133
$ws.Cells.Item(147, 5).Value = 103.775797

$ws.Cells.Item(148, 1).Value = "DT6"
$ws.Cells.Item(148, 2).Value = "King Albert Park"
$ws.Cells.Item(148, 4).Value = 1.335628
$ws.Cells.Item(148, 5).Value = 103.783983

$ws.Cells.Item(149, 1).Value = "DT7"
$ws.Cells.Item(149, 2).Value = "Sixth Avenue"
$ws.Cells.Item(149, 4).Value = 1.330714
$ws.Cells.Item(149, 5).Value = 103.797633

$ws.Cells.Item(150, 1).Value = "DT8"
$ws.Cells.Item(150, 2).Value = "Tan Kah Kee"
$ws.Cells.Item(150, 4).Value = 1.326039
$ws.Cells.Item(150, 5).Value = 103.807169

$ws.Cells.Item(151, 1).Value = "DT9"
$ws.Cells.Item(151, 2).Value = "Botanic Gardens"
$ws.Cells.Item(151, 3).Value = "植物园站"
$ws.Cells.Item(151, 4).Value = 1.322509
$ws.Cells.Item(151, 5).Value = 103.815376

$ws.Cells.Item(152, 1).Value = "DT10"
$ws.Cells.Item(152, 2).Value = "Stevens"
$ws.Cells.Item(152, 4).Value = 1.320069
$ws.Cells.Item(152, 5).Value = 103.825997

$ws.Cells.Item(153, 1).Value = "DT11"
$ws.Cells.Item(153, 2).Value = "Newton"
$ws.Cells.Item(153, 3).Value = "纽顿站"
$ws.Cells.Item(153, 4).Value = 1.312487
$ws.Cells.Item(153, 5).Value = 103.837924
$ws.Cells.Item(153, 6).Value = 228234
$ws.Cells.Item(153, 7).Value = "49 Scotts Rd"

$ws.Cells.Item(154, 1).Value = "DT12"
$ws.Cells.Item(154, 2).Value = "Little India"
$ws.Cells.Item(154, 3).Value = "小印度站"
$ws.Cells.Item(154, 4).Value = 1.307228
$ws.Cells.Item(154, 5).Value = 103.849847

$ws.Cells.Item(155, 1).Value = "DT13"
$ws.Cells.Item(155, 2).Value = "Rochor"
$ws.Cells.Item(155, 4).Value = 1.303764
$ws.Cells.Item(155, 5).Value = 103.852581

Write-Output "applied DT line 2 station updates"
